$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36 (shifts existing rows 36..105 down to 37..106)
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with its data
$ws.Cells.Item(36, 1).Value = 7
$ws.Cells.Item(36, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(36, 3).Value = "Ñuble"
$ws.Cells.Item(36, 4).Value = 44952
$ws.Cells.Item(36, 5).Value = 16
$ws.Cells.Item(36, 6).Value = 100112022
$ws.Cells.Item(36, 7).Value = "Arveja Verde"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 50
$ws.Cells.Item(36, 11).Value = 21000
$ws.Cells.Item(36, 12).Value = 21000
$ws.Cells.Item(36, 13).Value = 21000
$ws.Cells.Item(36, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(36, 15).Value = "Región de Ñuble"
$ws.Cells.Item(36, 16).Value = 840
$ws.Cells.Item(36, 17).Value = 25
$ws.Cells.Item(36, 18).Value = "Hortaliza"
